$wb = $excel.ActiveWorkbook

# --- "Measures" sheet -------------------------------------------------
$wsMeasures = $wb.Worksheets.Item("Measures")

# B2: simplify the DAX measure expression onto a single line
$wsMeasures.Range("B2").Value = "`nSUM('Sales SalesOrderDetail'[DerivedLineTotal])"

# E2: reword the measure description
$wsMeasures.Range("E2").Value = "This calculation is the sum of the line total derived from the SalesOrderDetail table. The line total value is the unit price of a particular item multiplied by the quantity of that item in a given sales order. This provides an overall total of the sales order."

# --- "Source Information" sheet ---------------------------------------
$wsSource = $wb.Worksheets.Item("Source Information")

# I2: reword the modification description
$wsSource.Range("I2").Value = "1. This round offs all the line item totals to the nearest 2 decimal places.`n`n"
